# Fix "Recorded By" column (G) ordering: swap the first and last
# comma-separated entries in each cell that contains more than one
# entry (cells with a single entry are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -gt 1) {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
        $newVal = [string]::Join(", ", $parts)
        $cell.Value2 = $newVal
    }
}
